$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 0
$ws.Range("H19").Value = 10382.19
$ws.Range("I19").Value = 852.3333
$ws.Range("J19").Value = 17529.584
$ws.Range("K19").Value = 852.3333
$ws.Range("L19").Value = 17529.584
$ws.Range("M19").Value = -677.3333
$ws.Range("N19").Value = -17879.584
$ws.Range("H100").Value = 2543.2693
$ws.Range("I100").Value = 2285.4167
$ws.Range("J100").Value = 2764.2856
$ws.Range("K100").Value = 2285.4167
$ws.Range("L100").Value = 2764.2856
$ws.Range("M100").Value = -1744.4167
$ws.Range("N100").Value = -3846.2856
$ws.Range("H103").Value = 125531.125
$ws.Range("I103").Value = 200409.8
$ws.Range("J103").Value = 733.3333
$ws.Range("K103").Value = 601229.3999999999
$ws.Range("L103").Value = 2199.9999
$ws.Range("M103").Value = -600643.3999999999
$ws.Range("N103").Value = -3371.9999
$ws.Range("H129").Value = 1050.3726
$ws.Range("I129").Value = 559.5
$ws.Range("J129").Value = 1236.1082
$ws.Range("K129").Value = 1678.5
$ws.Range("L129").Value = 3708.3246
$ws.Range("M129").Value = 3321.5
$ws.Range("N129").Value = -13708.3246
$ws.Range("H137").Value = 5606.909
$ws.Range("I137").Value = 5606.909
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 16820.727
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -14270.727
$ws.Range("H138").Value = 131511.03
$ws.Range("I138").Value = 2169.9
$ws.Range("J138").Value = 171930.12
$ws.Range("K138").Value = 6509.700000000001
$ws.Range("L138").Value = 515790.36
$ws.Range("M138").Value = -1369.700000000001
$ws.Range("N138").Value = -526070.36

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1450
$ws.Range("I25").Value = 1450
$ws.Range("K25").Value = 1450
$ws.Range("M25").Value = -1048
$ws.Range("H32").Value = 404527.5
$ws.Range("I32").Value = 448603.1
$ws.Range("K32").Value = 448603.1
$ws.Range("M32").Value = -448316.1
$ws.Range("H36").Value = 60021.75
$ws.Range("I36").Value = 30000
$ws.Range("K36").Value = 30000
$ws.Range("M36").Value = -29654
$ws.Range("H63").Value = 4787.3887
$ws.Range("I63").Value = 3454
$ws.Range("J63").Value = 6120.778
$ws.Range("K63").Value = 3454
$ws.Range("L63").Value = 6120.778
$ws.Range("M63").Value = -2768
$ws.Range("N63").Value = -7492.778
$ws.Range("H66").Value = 4787.3887
$ws.Range("I66").Value = 3454
$ws.Range("J66").Value = 6120.778
$ws.Range("K66").Value = 17270
$ws.Range("L66").Value = 30603.89
$ws.Range("M66").Value = -13838
$ws.Range("N66").Value = -37467.89
$ws.Range("H107").Value = 29800
$ws.Range("J107").Value = 29800
$ws.Range("L107").Value = 29800
$ws.Range("N107").Value = -37480
$ws.Range("H109").Value = 42000
$ws.Range("J109").Value = 42000
$ws.Range("L109").Value = 42000
$ws.Range("N109").Value = -44774

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1164.2222
$ws.Range("I37").Value = 1122.25
$ws.Range("J37").Value = 1500
$ws.Range("K37").Value = 1122.25
$ws.Range("L37").Value = 1500
$ws.Range("M37").Value = -985.25
$ws.Range("N37").Value = -1774
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").ClearContents()
$ws.Range("N38").Value = 0
$ws.Range("H86").Value = 250003000
$ws.Range("J86").Value = 5000
$ws.Range("L86").Value = 5000
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 250003000
$ws.Range("J89").Value = 5000
$ws.Range("L89").Value = 25000
$ws.Range("N89").Value = -36232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 14766.223
$ws.Range("J15").Value = 19649.334
$ws.Range("L15").Value = 19649.334
$ws.Range("N15").Value = -19989.334
$ws.Range("H31").Value = 2930.0193
$ws.Range("I31").Value = 990
$ws.Range("K31").Value = 990
$ws.Range("M31").Value = -695
$ws.Range("H34").Value = 2930.0193
$ws.Range("I34").Value = 990
$ws.Range("K34").Value = 990
$ws.Range("M34").Value = -788
$ws.Range("H132").Value = 15153852
$ws.Range("I132").Value = 979.8
$ws.Range("J132").Value = 27781244
$ws.Range("K132").Value = 2939.4
$ws.Range("L132").Value = 83343732
$ws.Range("M132").Value = -409.3999999999996
$ws.Range("N132").Value = -83348792

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 95
$ws.Range("I2").Value = 18.75
$ws.Range("J2").Value = 133.125
$ws.Range("K2").Value = 112.5
$ws.Range("L2").Value = 798.75
$ws.Range("M2").Value = 0.5
$ws.Range("N2").Value = -1024.75
$ws.Range("H21").Value = 2138.5
$ws.Range("I21").Value = 1518.3334
$ws.Range("J21").Value = 3999
$ws.Range("K21").Value = 4555.0002
$ws.Range("L21").Value = 11997
$ws.Range("M21").Value = -4382.0002
$ws.Range("N21").Value = -12343
$ws.Range("H97").Value = 29001
$ws.Range("I97").Value = 600
$ws.Range("J97").Value = 100003.5
$ws.Range("K97").Value = 1800
$ws.Range("L97").Value = 300010.5
$ws.Range("M97").Value = -1304
$ws.Range("N97").Value = -301002.5
$ws.Range("H140").Value = 1617.2
$ws.Range("J140").Value = 3658
$ws.Range("L140").Value = 10974
$ws.Range("N140").Value = -21334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 46557.332
$ws.Range("J22").Value = 46557.332
$ws.Range("L22").Value = 46557.332
$ws.Range("N22").Value = -47615.332
$ws.Range("H24").Value = 5001262.5
$ws.Range("J24").Value = 1683.3334
$ws.Range("L24").Value = 1683.3334
$ws.Range("N24").Value = -2029.3334
$ws.Range("H28").Value = 28999
$ws.Range("J28").Value = 28999
$ws.Range("L28").Value = 28999
$ws.Range("N28").Value = -29383
$ws.Range("H29").Value = 8500
$ws.Range("J29").Value = 8500
$ws.Range("L29").Value = 8500
$ws.Range("N29").Value = -9080
$ws.Range("H31").Value = 2200
$ws.Range("I31").Value = 2200
$ws.Range("K31").Value = 2200
$ws.Range("M31").Value = -1908
$ws.Range("H37").Value = 2200
$ws.Range("I37").Value = 2200
$ws.Range("K37").Value = 2200
$ws.Range("M37").Value = -1923
$ws.Range("H41").Value = 2669.2
$ws.Range("J41").Value = 850
$ws.Range("L41").Value = 850
$ws.Range("N41").Value = -1560
$ws.Range("H93").Value = 42998.824
$ws.Range("J93").Value = 42998.824
$ws.Range("L93").Value = 42998.824
$ws.Range("N93").Value = -46742.824
$ws.Range("H97").Value = 2241.5
$ws.Range("I97").Value = 2176.25
$ws.Range("K97").Value = 2176.25
$ws.Range("M97").Value = -1680.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 56671.668
$ws.Range("I4").Value = 33333.332
$ws.Range("J4").Value = 80010
$ws.Range("K4").Value = 33333.332
$ws.Range("L4").Value = 80010
$ws.Range("M4").Value = -33220.332
$ws.Range("N4").Value = -80236
$ws.Range("H22").Value = 10610
$ws.Range("I22").Value = 920.2
$ws.Range("J22").Value = 18684.834
$ws.Range("K22").Value = 920.2
$ws.Range("L22").Value = 18684.834
$ws.Range("M22").Value = -625.2
$ws.Range("N22").Value = -19274.834
$ws.Range("H26").Value = 54755
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9705
$ws.Range("H27").Value = 10610
$ws.Range("I27").Value = 920.2
$ws.Range("J27").Value = 18684.834
$ws.Range("K27").Value = 920.2
$ws.Range("L27").Value = 18684.834
$ws.Range("M27").Value = -813.2
$ws.Range("N27").Value = -18898.834
$ws.Range("H28").Value = 56671.668
$ws.Range("I28").Value = 33333.332
$ws.Range("J28").Value = 80010
$ws.Range("K28").Value = 33333.332
$ws.Range("L28").Value = 80010
$ws.Range("M28").Value = -33101.332
$ws.Range("N28").Value = -80474
$ws.Range("H29").Value = 33354852
$ws.Range("I29").Value = 55555
$ws.Range("J29").Value = 50004500
$ws.Range("K29").Value = 55555
$ws.Range("L29").Value = 50004500
$ws.Range("M29").Value = -55260
$ws.Range("N29").Value = -50005090
$ws.Range("H37").Value = 56671.668
$ws.Range("I37").Value = 33333.332
$ws.Range("J37").Value = 80010
$ws.Range("K37").Value = 33333.332
$ws.Range("L37").Value = 80010
$ws.Range("M37").Value = -33226.332
$ws.Range("N37").Value = -80224
$ws.Range("H43").Value = 34148.855
$ws.Range("J43").Value = 34148.855
$ws.Range("L43").Value = 34148.855
$ws.Range("N43").Value = -34534.855
$ws.Range("H132").Value = 4155.237
$ws.Range("I132").Value = 3524.238
$ws.Range("J132").Value = 4934.706
$ws.Range("K132").Value = 10572.714
$ws.Range("L132").Value = 14804.118
$ws.Range("M132").Value = -8042.714
$ws.Range("N132").Value = -19864.118
$ws.Range("H136").Value = 5052219
$ws.Range("I136").Value = 1759.3636
$ws.Range("J136").Value = 15153138
$ws.Range("K136").Value = 5278.0908
$ws.Range("L136").Value = 45459414
$ws.Range("M136").Value = -2728.0908
$ws.Range("N136").Value = -45464514

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1016501.5
$ws.Range("I4").Value = 3000000
$ws.Range("J4").Value = 24752.25
$ws.Range("K4").Value = 3000000
$ws.Range("L4").Value = 24752.25
$ws.Range("M4").Value = -2999887
$ws.Range("N4").Value = -24978.25
$ws.Range("H6").Value = 23803.666
$ws.Range("I6").Value = 205
$ws.Range("K6").Value = 205
$ws.Range("M6").Value = -90
$ws.Range("H31").Value = 56611.4
$ws.Range("J31").Value = 56611.4
$ws.Range("L31").Value = 56611.4
$ws.Range("N31").Value = -57307.4
